$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 8.660232485948974;  G = 17.45944343273191 }
    3 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    4 = @{ B = 0.6753301551942219; C = 0.3127903958511391; D = 0.1575252929769615; E = 0.496779210170732;  G = 1.642425054193055 }
    5 = @{ B = 1.459612070389937;  C = 10.29869402782916;  D = 0.8054896365839992; E = 8.660232485948974;  G = 21.22402822075207 }
    6 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 3.900430680208489;  E = 0.496779210170732;  G = 9.295990156953671 }
    7 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 0.1575252929769615; E = 0.496779210170732;  G = 2.997429241610044 }
    8 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 3.645393585217082 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
